$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting from column I into the new column J (rows 3-12) so the
# new column inherits the same styles (borders, number format, etc.) as the
# existing "2019" column.
$ws.Range("I3:I12").Copy($ws.Range("J3:J12"))

# Row 4 header: new year column
$ws.Range("J4").Value2 = 2020

# Row 5 values
$ws.Range("J5").Value2 = 253.27664777870578

# Row 6 stays blank (section header row) - already blank after the copy.

# Row 7 values
$ws.Range("J7").Value2 = 93.236077839070575

# Row 8 values
$ws.Range("J8").Value2 = 160

# Row 9 stays blank (section header row) - already blank after the copy.

# Row 10 values
$ws.Range("J10").Value2 = 69

# Row 11 values
$ws.Range("J11").Value2 = 48.5

# Row 12 values
$ws.Range("J12").Value2 = 22.8

# Update the active selection to reflect where the user left off editing.
$ws.Range("J3").Select()
